$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the header row (row 1), pushing existing data down.
$ws.Rows.Item(2).Resize(2).EntireRow.Insert() | Out-Null

# New row 3: mean_temperature / mean / LST (deg C)   -- set first so "mean_temperature"
# gets interned into the shared-string table ahead of "rainfall".
$ws.Cells.Item(3, 1).Value = "mean_temperature"
$ws.Cells.Item(3, 2).Value = "mean"
$ws.Cells.Item(3, 3).Value = "LST (" + [char]0x00B0 + "C)"

# New row 2: rainfall / sum / Rain (mm)
$ws.Cells.Item(2, 1).Value = "rainfall"
$ws.Cells.Item(2, 2).Value = "sum"
$ws.Cells.Item(2, 3).Value = "Rain (mm)"

# Fix the totprec row (now row 4) - should be sum / Rain (mm), not mean / LST
$ws.Cells.Item(4, 2).Value = "sum"
$ws.Cells.Item(4, 3).Value = "Rain (mm)"

# Column A width update (closest achievable value given engine rounding)
$ws.Columns.Item(1).ColumnWidth = 17.6

# Selection update
$ws.Range("H12").Select() | Out-Null
